# Insert a new row at row 105 (shifting existing rows 105..220 down to 106..221)
# and populate the new row 105 with the new data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("105:105").Insert()

$ws.Range("A105").Value = 3
$ws.Range("B105").Value = "Femacal de La Calera"
$ws.Range("C105").Value = "Coquimbo"
$ws.Range("D105").Value = 44467
$ws.Range("E105").Value = 5
$ws.Range("F105").Value = 100112003
$ws.Range("G105").Value = "Ajo"
$ws.Range("H105").Value = "Chino"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 78
$ws.Range("K105").Value = 16000
$ws.Range("L105").Value = 16500
$ws.Range("M105").Value = 16244
$ws.Range("N105").Value = '$/caja 10 kilos'
$ws.Range("O105").Value = "China"
$ws.Range("P105").Value = 1624
$ws.Range("Q105").Value = 10
$ws.Range("R105").Value = "Hortaliza"

# Make sure D105 keeps the same date style as the other D-column cells (s="2")
$ws.Range("D105").NumberFormat = $ws.Range("D106").NumberFormat
